$wb = $excel.ActiveWorkbook

# Column F ("想去人数" / "want-to-go count") updates, applied identically to
# the "展览" (sheet1) and "全部类型" (sheet4) sheets.
$updates = @{
    "F2"  = 211
    "F3"  = 437
    "F4"  = 12866
    "F6"  = 190
    "F16" = 396
    "F17" = 5489
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
